$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated measurements (columns A-C) and normalized "busque" label (column F)
# for rows 2-11, per the 4 de julho de 2025 dataset refresh.

$ws.Range("A2").Value = 37.6529892109601
$ws.Range("B2").Value = 711
$ws.Range("C2").Value = 619
$ws.Range("F2").Value = "busque"

$ws.Range("A3").Value = 87.23174527890686
$ws.Range("B3").Value = 424
$ws.Range("C3").Value = 689
$ws.Range("F3").Value = "busque"

$ws.Range("A4").Value = 156.2591062507433
$ws.Range("B4").Value = 0
$ws.Range("C4").Value = 444
$ws.Range("F4").Value = "busque"

$ws.Range("A5").Value = 158.5072909037546
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = 590
$ws.Range("F5").Value = "busque"

$ws.Range("A6").Value = 155.6952038123114
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = 462
$ws.Range("F6").Value = "busque"

$ws.Range("A7").Value = 123.6853652880928
$ws.Range("B7").Value = 192
$ws.Range("C7").Value = 625
$ws.Range("F7").Value = "busque"

$ws.Range("A8").Value = 33.15462804192198
$ws.Range("B8").Value = 707
$ws.Range("C8").Value = 590
$ws.Range("F8").Value = "busque"

$ws.Range("A9").Value = 163.3186278728956
$ws.Range("B9").Value = 4
$ws.Range("C9").Value = 704
$ws.Range("F9").Value = "busque"

$ws.Range("A10").Value = 56.29359121085295
$ws.Range("B10").Value = 524
$ws.Range("C10").Value = 531
$ws.Range("F10").Value = "busque"

$ws.Range("A11").Value = 105.0520715571112
$ws.Range("B11").Value = 307
$ws.Range("C11").Value = 663
$ws.Range("F11").Value = "busque"
